$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold, border, centered) from H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows
$data = @(
    @(1, 4),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 2),
    @(1, 6),
    @(7, 9),
    @(1, 6),
    @(1, 4),
    @(1, 3),
    @(5, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}

$excel.CutCopyMode = $false
